$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new record was logged for 2026/02/12 (Thursday) which needs to be
# inserted right after the existing 2026/02/11 rows (i.e. before the
# 2026/12/29 block), shifting every subsequent row down by one.
$ws.Rows(798).Insert()

# Column A stores the date as plain text (e.g. "2026/12/29"), not a real
# Excel date. Force a Text number format before writing the value so
# Excel's automation layer doesn't auto-convert the slash-separated
# string into a date serial number, then clear the formatting again so
# the cell ends up with the same (default) style as all of its neighbors.
$newRow = 798
$ws.Range("A" + $newRow).NumberFormat = "@"
$ws.Range("A" + $newRow).Value = "2026/02/12"
$ws.Range("A" + $newRow).ClearFormats()

$ws.Range("B" + $newRow).Value = "木"
$ws.Range("C" + $newRow).Value = 2
$ws.Range("D" + $newRow).Value = 201
